$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 203.5
$ws.Range("I9").Value = 75
$ws.Range("K9").Value = 75
$ws.Range("M9").Value = 94

$ws.Range("H11").Value = 487.52942
$ws.Range("I11").Value = 487.52942
$ws.Range("K11").Value = 487.52942
$ws.Range("M11").Value = -347.52942

$ws.Range("H28").Value = 1624.7693
$ws.Range("J28").Value = 2065.1428
$ws.Range("L28").Value = 2065.1428
$ws.Range("N28").Value = -3035.1428

$ws.Range("H41").Value = 1049.5
$ws.Range("J41").Value = 2205
$ws.Range("L41").Value = 2205
$ws.Range("N41").Value = -3085

$ws.Range("H92").Value = 111111656
$ws.Range("I92").Value = 142857760
$ws.Range("J92").Value = 320
$ws.Range("K92").Value = 142857760
$ws.Range("L92").Value = 320
$ws.Range("M92").Value = -142856512
$ws.Range("N92").Value = -2816

$ws.Range("H106").Value = 3139.4546
$ws.Range("I106").Value = 3069.5715
$ws.Range("K106").Value = 3069.5715
$ws.Range("M106").Value = -2438.5715

$ws.Range("H118").Value = 454.85
$ws.Range("I118").Value = 467.73685
$ws.Range("J118").Value = 210
$ws.Range("K118").Value = 1403.21055
$ws.Range("L118").Value = 630
$ws.Range("M118").Value = 253.78945
$ws.Range("N118").Value = -3944

$ws.Range("H132").Value = 3626.2205
$ws.Range("I132").Value = 1221.46
$ws.Range("K132").Value = 3664.38
$ws.Range("M132").Value = -1134.38

$ws.Range("H138").Value = 5447.86
$ws.Range("I138").Value = 1726.579
$ws.Range("J138").Value = 6320.753
$ws.Range("K138").Value = 5179.737
$ws.Range("L138").Value = 18962.259
$ws.Range("M138").Value = -39.73700000000008
$ws.Range("N138").Value = -29242.259

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11733.262
$ws.Range("I32").Value = 11152.776
$ws.Range("K32").Value = 11152.776
$ws.Range("M32").Value = -10865.776

$ws.Range("H74").Value = 13889788
$ws.Range("I74").Value = 14706720
$ws.Range("K74").Value = 14706720
$ws.Range("M74").Value = -14705846

$ws.Range("H77").Value = 13889788
$ws.Range("I77").Value = 14706720
$ws.Range("K77").Value = 73533600
$ws.Range("M77").Value = -73529232

$ws.Range("H88").Value = 3091.2
$ws.Range("I88").Value = 1000
$ws.Range("J88").Value = 3323.5557
$ws.Range("K88").Value = 1000
$ws.Range("L88").Value = 3323.5557
$ws.Range("M88").Value = -594
$ws.Range("N88").Value = -4135.5557

$ws.Range("H91").Value = 3091.2
$ws.Range("I91").Value = 1000
$ws.Range("J91").Value = 3323.5557
$ws.Range("K91").Value = 1000
$ws.Range("L91").Value = 3323.5557
$ws.Range("M91").Value = 404
$ws.Range("N91").Value = -6131.5557

$ws.Range("H110").Value = 1630
$ws.Range("I110").Value = 1372.2
$ws.Range("K110").Value = 1372.2
$ws.Range("M110").Value = 672.8

$ws.Range("H132").Value = 14101
$ws.Range("I132").Value = 16587.158
$ws.Range("J132").Value = 3603.889
$ws.Range("K132").Value = 49761.474
$ws.Range("L132").Value = 10811.667
$ws.Range("M132").Value = -47231.474
$ws.Range("N132").Value = -15871.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2603.3928
$ws.Range("I86").Value = 1691.8636
$ws.Range("K86").Value = 1691.8636
$ws.Range("M86").Value = -568.8635999999999

$ws.Range("H89").Value = 2603.3928
$ws.Range("I89").Value = 1691.8636
$ws.Range("K89").Value = 8459.317999999999
$ws.Range("M89").Value = -2843.317999999999

$ws.Range("H96").Value = 37824.668
$ws.Range("I96").Value = 27750
$ws.Range("J96").Value = 57974
$ws.Range("K96").Value = 27750
$ws.Range("L96").Value = 57974
$ws.Range("M96").Value = -25004
$ws.Range("N96").Value = -63466

$ws.Range("H105").Value = 4722
$ws.Range("I105").Value = 3963
$ws.Range("K105").Value = 3963
$ws.Range("M105").Value = -2216

$ws.Range("H107").Value = 2525.3823
$ws.Range("I107").Value = 2195.9614
$ws.Range("J107").Value = 3596
$ws.Range("K107").Value = 2195.9614
$ws.Range("L107").Value = 3596
$ws.Range("M107").Value = -275.9614000000001
$ws.Range("N107").Value = -7436

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 192
$ws.Range("I7").Value = 168.71428
$ws.Range("J7").Value = 215.28572
$ws.Range("K7").Value = 168.71428
$ws.Range("L7").Value = 215.28572
$ws.Range("M7").Value = -55.71428
$ws.Range("N7").Value = -441.28572

$ws.Range("H31").Value = 21742582
$ws.Range("I31").Value = 25643976
$ws.Range("K31").Value = 25643976
$ws.Range("M31").Value = -25643681

$ws.Range("H34").Value = 21742582
$ws.Range("I34").Value = 25643976
$ws.Range("K34").Value = 25643976
$ws.Range("M34").Value = -25643774

$ws.Range("H107").Value = 865.93335
$ws.Range("I107").Value = 510.1111
$ws.Range("K107").Value = 510.1111
$ws.Range("M107").Value = 1409.8889

$ws.Range("H132").Value = 43012668
$ws.Range("J132").Value = 2999.1428
$ws.Range("L132").Value = 8997.428400000001
$ws.Range("N132").Value = -14057.4284

$ws.Range("H133").Value = 88325.5
$ws.Range("J133").Value = 88325.5
$ws.Range("L133").Value = 88325.5
$ws.Range("N133").Value = -93385.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 75.09090999999999
$ws.Range("I2").Value = 60
$ws.Range("J2").Value = 80.75
$ws.Range("K2").Value = 360
$ws.Range("L2").Value = 484.5
$ws.Range("M2").Value = -247
$ws.Range("N2").Value = -710.5

$ws.Range("H107").Value = 2277.3333
$ws.Range("I107").Value = 2749.5
$ws.Range("J107").Value = 2142.4285
$ws.Range("K107").Value = 8248.5
$ws.Range("L107").Value = 6427.2855
$ws.Range("M107").Value = -6328.5
$ws.Range("N107").Value = -10267.2855

$ws.Range("H131").Value = 8335884.5
$ws.Range("J131").Value = 10103873
$ws.Range("L131").Value = 30311619
$ws.Range("N131").Value = -30321699

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

$ws.Range("H80").Value = 93137.08
$ws.Range("J80").Value = 28481.166
$ws.Range("L80").Value = 28481.166
$ws.Range("N80").Value = -30477.166

$ws.Range("H83").Value = 93137.08
$ws.Range("J83").Value = 28481.166
$ws.Range("L83").Value = 142405.83
$ws.Range("N83").Value = -152389.83

$ws.Range("H97").Value = 984.129
$ws.Range("I97").Value = 754.4783
$ws.Range("J97").Value = 1644.375
$ws.Range("K97").Value = 754.4783
$ws.Range("L97").Value = 1644.375
$ws.Range("M97").Value = -258.4783
$ws.Range("N97").Value = -2636.375

$ws.Range("H102").Value = 15940627
$ws.Range("I102").Value = 19618438
$ws.Range("K102").Value = 19618438
$ws.Range("M102").Value = -19616816

$ws.Range("H122").Value = 195803.56
$ws.Range("I122").Value = 265309.75
$ws.Range("J122").Value = 7143.9287
$ws.Range("K122").Value = 795929.25
$ws.Range("L122").Value = 21431.7861
$ws.Range("M122").Value = -793479.25
$ws.Range("N122").Value = -26331.7861

$ws.Range("H132").Value = 86805.53999999999
$ws.Range("I132").Value = 120430.586
$ws.Range("K132").Value = 361291.758
$ws.Range("M132").Value = -358761.758

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2231.4
$ws.Range("I61").Value = 1300.2858
$ws.Range("K61").Value = 1300.2858
$ws.Range("M61").Value = -1098.2858

$ws.Range("H82").Value = 2483.8823
$ws.Range("I82").Value = 570.6667
$ws.Range("K82").Value = 570.6667
$ws.Range("M82").Value = -209.6667

$ws.Range("H85").Value = 2483.8823
$ws.Range("I85").Value = 570.6667
$ws.Range("K85").Value = 570.6667
$ws.Range("M85").Value = 677.3333

$ws.Range("H96").Value = 40000
$ws.Range("J96").Value = 40000
$ws.Range("L96").Value = 40000
$ws.Range("N96").Value = -45492

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H113").Value = 2231.4
$ws.Range("I113").Value = 1300.2858
$ws.Range("K113").Value = 1300.2858
$ws.Range("M113").Value = 869.7141999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1828.4166
$ws.Range("J96").Value = 2098.25
$ws.Range("L96").Value = 2098.25
$ws.Range("N96").Value = -4844.25

$ws.Range("H136").Value = 3741.457
$ws.Range("I136").Value = 2933.9092
$ws.Range("J136").Value = 5108.077
$ws.Range("K136").Value = 8801.7276
$ws.Range("L136").Value = 15324.231
$ws.Range("M136").Value = -6251.7276
$ws.Range("N136").Value = -20424.231
